$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.321.19'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '2.228.14'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''244.50'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").Value = '''0.628'
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("D7").Value = '''73.93'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").Value = '''42.55'
$ws.Range("E10").Value = '  +4.74%  '
$ws.Range("D11").Value = '''0.0972'
$ws.Range("E11").Value = '  +4.22%  '
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").Value = '''14.40'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = '''0.852'
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '2.226.78'
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.0000117'
$ws.Range("E17").Value = '  +20.69%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '42.130.88'
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("E19").Value = '  +2.60%  '
$ws.Range("D20").Value = '''72.13'
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").Value = '''9.92'
$ws.Range("E21").Value = '  +38.16%  '
$ws.Range("D22").Value = '''231.12'
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("D23").Value = '''2.14'
$ws.Range("E23").Value = '  -3.69%  '
$ws.Range("D24").Value = '''11.86'
$ws.Range("E24").Value = '  +8.15%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("E28").Value = '  +3.33%  '
$ws.Range("D29").Value = '''167.18'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D31").Value = '''5.74'
$ws.Range("E31").Value = '  +19.64%  '
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("E33").Value = '  +0.23%  '
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("D35").Value = '''29.47'
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("D36").Value = '''4.41'
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("E37").Value = '  +2.96%  '
$ws.Range("D38").Value = '''13.00'
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("E39").Value = '  +1.34%  '
$ws.Range("D40").Value = '''5.61'
$ws.Range("E40").Value = '  -1.59%  '
$ws.Range("D41").Value = '''62.78'
$ws.Range("E41").Value = '  +5.56%  '
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("E43").Value = '  +2.74%  '
$ws.Range("D44").Value = '''105.37'
$ws.Range("E44").Value = '  -4.55%  '
$ws.Range("E45").Value = '  +3.48%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("E47").Value = '  +7.32%  '
$ws.Range("E48").Value = '  +1.44%  '
$ws.Range("D49").Value = '''1.18'
$ws.Range("E49").Value = '  +2.64%  '
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("D51").Value = '''4.05'
$ws.Range("E51").Value = '  +0.43%  '
